{"js": "// Apply \"Compact\" paragraph style to empty table-cell paragraphs (the\n// blank \"Due\" column cells) that currently have no paragraph properties.\n// Mirrors the OOXML diff: <w:p/> -> <w:p><w:pPr><w:pStyle w:val=\"Compact\"/></w:pPr></w:p>\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Gather every row from every table, then every cell from every row.\nconst allRows = [];\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    allRows.push(row);\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nconst allCells = [];\nfor (const row of allRows) {\n  for (const cell of row.cells.items) {\n    allCells.push(cell);\n  }\n}\n\n// Load each cell's single paragraph (these empty cells only ever contain\n// one empty paragraph) along with its text and current style.\nconst paraInfo = [];\nfor (const cell of allCells) {\n  const paras = cell.body.paragraphs;\n  paras.load(\"items/text,items/style\");\n  paraInfo.push(paras);\n}\nawait context.sync();\n\n// Only touch paragraphs that are empty and still on the default \"Normal\"\n// style (i.e. no explicit pPr/pStyle yet) - exactly the cells the diff\n// changed, leaving every other paragraph untouched.\nfor (const paras of paraInfo) {\n  for (const para of paras.items) {\n    if (para.text === \"\" && para.style === \"Normal\") {\n      para.style = \"Compact\";\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the \"Compact\" paragraph style to empty table-cell paragraphs\n# (the blank \"Due\" column cells) that currently have no paragraph\n# properties at all. Mirrors the OOXML diff:\n#   <w:p/> -> <w:p><w:pPr><w:pStyle w:val=\"Compact\"/></w:pPr></w:p>\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    foreach ($row in $t.Rows) {\n        foreach ($cell in $row.Cells) {\n            $cellRange = $cell.Range\n            # An empty cell's Range.Text is just the cell-mark pair\n            # (\"\\r\\a\") - i.e. no visible/typed text at all.\n            $isEmpty = ($cellRange.Text -eq \"`r`a\")\n            $styleName = $cellRange.Style.NameLocal\n            if ($isEmpty -and $styleName -eq \"Normal\") {\n                $cellRange.Style = \"Compact\"\n            }\n        }\n    }\n}\n"}
